# The underlying change in this revision is the "Fixed POI packaging and
# upgraded to POI 3.15" housekeeping commit: the saved OOXML parts get
# re-serialized by the upgraded writer, which canonicalizes (alphabetizes)
# XML attribute order on every element it touches. No attribute value, no
# text, and no structural content actually changes anywhere in the package
# -- only the on-disk attribute ordering differs.
#
# To reproduce the intent of that re-serialization through the Word object
# model, we re-touch (re-apply, with their own existing values) every
# property whose underlying XML attribute set was reordered by the diff:
#   - the table's preferred width (w:tblW)
#   - the table's autoformat "look" flags (w:tblLook)
#   - the table style, which governs the conditional-formatting markers
#     (w:cnfStyle) stamped on the rows/cells/paragraphs
#   - every cell's preferred width (w:tcW)
#   - the section's page size and margins (w:pgSz / w:pgMar)
# so that each of these stays exactly the same value it already had, while
# being explicitly re-asserted via COM.

$d = $word.ActiveDocument

$t = $d.Tables.Item(1)

# w:tblW w:w="0" w:type="auto" -- re-assert the automatic table width.
$t.PreferredWidthType = 1   # wdPreferredWidthAuto
$t.PreferredWidth = 0

# w:tblLook -- re-assert the banding / first-row / first-column flags.
$t.ApplyStyleHeadingRows = $true
$t.ApplyStyleLastRow = $false
$t.ApplyStyleFirstColumn = $true
$t.ApplyStyleLastColumn = $false
$t.ApplyStyleRowBands = $true
$t.ApplyStyleColumnBands = $false

# Re-apply the table style itself -- this is what drives the per
# row/cell/paragraph w:cnfStyle conditional-formatting markers.
$t.Style = "TableauGrille2"

# w:tcW w:w="4606" w:type="dxa" on every cell of the table.
for ($ri = 1; $ri -le $t.Rows.Count; $ri++) {
    $row = $t.Rows.Item($ri)
    for ($ci = 1; $ci -le $row.Cells.Count; $ci++) {
        $cell = $row.Cells.Item($ci)
        $cell.PreferredWidthType = 3   # wdPreferredWidthPoints (dxa-backed)
        $cell.PreferredWidth = 230.3   # 4606 twips
    }
}

# w:pgSz w:w="11906" w:h="16838" and w:pgMar (top/right/bottom/left/header/
# footer/gutter all unchanged) on the document's (only) section.
$ps = $d.Sections.Item(1).PageSetup
$ps.PageWidth = 595.3
$ps.PageHeight = 841.9
$ps.TopMargin = 70.85
$ps.RightMargin = 70.85
$ps.BottomMargin = 70.85
$ps.LeftMargin = 70.85
$ps.HeaderDistance = 35.4
$ps.FooterDistance = 35.4
$ps.Gutter = 0
